$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

# Row 10
$ws.Range("E10").Value = 520

# Row 11
$ws.Range("E11").Value = 336
$ws.Range("G11").Value = 67
$ws.Range("H11").Value = 253

# Row 12
$ws.Range("E12").Value = 509
$ws.Range("G12").Value = 85
$ws.Range("H12").Value = 361

# Row 13
$ws.Range("E13").Value = 127

# Row 30
$ws.Range("E30").Value = 202

# Row 49
$ws.Range("E49").Value = 282

# Row 51
$ws.Range("E51").Value = 228
